# Update the tech table on Sheet2.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# --- Row 10 (Wifi) : fill in the comparison columns H:K ---
$ws.Range("H10").Value = "Weak"
$ws.Range("I10").Value = "Strong"
$ws.Range("J10").Value = "Strong"
$ws.Range("K10").Value = "High"

# --- Row 11 (Zigbee) : fill in the comparison columns H:K ---
$ws.Range("H11").Value = "Strong"
$ws.Range("I11").Value = "Strong"
$ws.Range("J11").Value = "Strong"
$ws.Range("K11").Value = "Average"

# --- Row 12 (NFC) : fill in columns D,E,G and H:K ---
$ws.Range("E12").Value = "yes"
$ws.Range("D12").Value = "250 kb/s"
$ws.Range("G12").Value = "up to 100m"
$ws.Range("H12").Value = "Strong"
$ws.Range("I12").Value = "Strong"
$ws.Range("J12").Value = "Strong"
$ws.Range("K12").Value = "low"

# --- Row 21 (Display driver) : fill in columns C:I ---
$ws.Range("C21").Value = "15v or 5v"
$ws.Range("D21").Value = ".05 uA / cm^2"
$ws.Range("E21").Value = "integrated driver"
$ws.Range("F21").Value = "unknown"
$ws.Range("G21").Value = "unknown"
$ws.Range("H21").Value = "varies"
$ws.Range("I21").Value = "varies"
# the longer "integrated driver" text wraps within its column, so the row
# grows to fit two lines of text
$ws.Rows.Item(21).RowHeight = 30

# --- Row 25 (Phone (Android)) : fill in columns C:F ---
$ws.Range("C25").Value = "Yes"
$ws.Range("D25").Value = "Yes"
$ws.Range("E25").Value = "yes"
$ws.Range("F25").Value = "No"

# --- Row 26 (Tablet (Andoid)) : correct columns D,E ---
$ws.Range("D26").Value = "Yes"
$ws.Range("E26").Value = "Yes"

# --- Row 28 (Windows phone (Windows 8)) : fill in columns C:F ---
$ws.Range("C28").Value = "Yes"
$ws.Range("D28").Value = "No"
$ws.Range("E28").Value = "No"
$ws.Range("F28").Value = "Yes"

# --- Row 29 (Windows table (Windows 8)) : fill in columns C:F ---
$ws.Range("C29").Value = "Yes"
$ws.Range("D29").Value = "No"
$ws.Range("E29").Value = "No"
$ws.Range("F29").Value = "No"

# --- Update the saved view: scroll position and active selection ---
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("G29").Select()
